$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- PasteSpecial paste-type constants ---
$xlPasteValues  = -4163
$xlPasteFormats = -4122

function Set-TextCell($row, $col, $text) {
    # Writing a literal value that looks like a date (e.g. "2022-08-03") would get
    # auto-converted into a real date serial number by a plain Value assignment.
    # To keep it as plain text (matching the source workbook, which stores these
    # values as text) we first write it as a formula evaluating to the literal
    # string, then flatten the formula down to a static value via copy /
    # paste-special values. This preserves the cell's existing style.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues)
}

function Set-RunTimeFormat($row) {
    # Column B holds a date+time serial number; give it the same custom
    # date/time number format already used throughout the "Run Time" column.
    $ws.Cells.Item(49, 2).Copy() | Out-Null
    $ws.Cells.Item($row, 2).PasteSpecial($xlPasteFormats)
}

function Set-RowValues($row, $a, $b, $c, $d, $e, $f, $g) {
    Set-TextCell $row 1 $a
    $ws.Cells.Item($row, 2).Value = $b
    Set-TextCell $row 3 $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g

    Set-RunTimeFormat $row
}

# ------------------------------------------------------------------
# Fix up existing row 50: re-apply the same "data row" style used by
# the rest of the table, and correct the slightly-off run time value
# in column B.
# ------------------------------------------------------------------
Set-RowValues 50 "2022-07-15" 44757.63793371528 "testcerti164" 51 51 0 1.31

# ------------------------------------------------------------------
# Append new history rows 51-57
# ------------------------------------------------------------------
Set-RowValues 51 "2022-08-03" 44776.69886240741 "165_scndcycle" 51 51 0 1.16
Set-RowValues 52 "2022-08-04" 44777.40108302084 "165finalrun"   51 51 0 1.17
Set-RowValues 53 "2022-08-22" 44795.68176572917 "166fstcycle"   51 51 0 1.31
Set-RowValues 54 "2022-08-23" 44796.9161703125  "166cyclescnd"  51 51 0 1.4
Set-RowValues 55 "2022-08-30" 44803.90564909722 "cert234"       51 47 4 1.15
Set-RowValues 56 "2022-08-30" 44803.90679255787 "cert345"       51 51 0 1.01
Set-RowValues 57 "2022-08-30" 44803.91025019104 "cert456"       51 51 0 0.99
